# Updates the cryptos price/volume table with the latest scraped values
# (GitHub Actions refresh run). A couple of rows also changed identity:
# row 45/46 (Aave <-> BabyDogeCoin) swapped rank order, and row 51
# (Algorand) fell out of the top 50 and was replaced by Mantle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column ("D") values that look like plain numbers (e.g. "1.000",
# "0.06440") must be forced to text first, otherwise Excel silently
# coerces them to numeric values and the significant trailing/leading
# zeros the source format relies on get lost.

$ws.Range("D2").Value = "26.515.74"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "1.686.69"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.83"
$ws.Range("E5").Value = "  +5.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5346"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06440"
$ws.Range("E9").Value = "  +3.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.41"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "1.694.46"
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.509"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5639"
$ws.Range("E14").Value = "  +5.54%  "
$ws.Range("D15").Value = "0.0₅8429"
$ws.Range("E15").Value = "  +5.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.36"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "26.551.12"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.837"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.17"
$ws.Range("E20").Value = "  +4.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.409"
$ws.Range("E22").Value = "  +5.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.68"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1276"
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.507"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.28"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.424"
$ws.Range("E28").Value = "  +4.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06166"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.281"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.570"
$ws.Range("E31").Value = "  +4.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.465"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.719"
$ws.Range("E33").Value = "  +6.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.022"
$ws.Range("E34").Value = "  +5.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.801"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5753"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01650"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.961"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8705"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").Value = "1.056.60"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.30"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "1.838.03"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈111"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.30"
$ws.Range("E46").Value = "  +5.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.182"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.120"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05207"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4231"
$ws.Range("E51").Value = "  -0.14%  "
